$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "2025-12-11 Thursday" "2025-12-12 Friday"

Replace-Text "900×2=1800" "552×4=2208"
Replace-Text "499×7=3493" "282×7=1974"
Replace-Text "354×5=1770" "200×8=1600"
Replace-Text "564×8=4512" "148×7=1036"
Replace-Text "177×9=1593" "920×2=1840"
Replace-Text "956×6=5736" "993×3=2979"
Replace-Text "890×9=8010" "271×5=1355"
Replace-Text "373×8=2984" "991×8=7928"
Replace-Text "331×6=1986" "467×7=3269"
Replace-Text "922×3=2766" "591×2=1182"
Replace-Text "292×2=584" "290×3=870"
Replace-Text "742×4=2968" "505×9=4545"
Replace-Text "526×9=4734" "939×3=2817"
Replace-Text "955×3=2865" "300×2=600"
Replace-Text "799×6=4794" "581×7=4067"
Replace-Text "107×8=856" "623×2=1246"
Replace-Text "477×9=4293" "569×8=4552"
Replace-Text "343×8=2744" "831×5=4155"
Replace-Text "958×9=8622" "924×3=2772"
Replace-Text "266×6=1596" "280×4=1120"
Replace-Text "112×6=672" "997×7=6979"
Replace-Text "783×4=3132" "764×4=3056"
Replace-Text "589×3=1767" "660×7=4620"
Replace-Text "809×6=4854" "812×8=6496"
Replace-Text "718×5=3590" "827×7=5789"
